{"js": "// Apply the report.docx content edits described by the diff:\n//  1. Rewrite the CHOLMOD ordering-method sentence (expand CHOLMOD_AMD /\n//     CHOLMOD_METIS into spelled-out English names, \"\u4e3b\u8981\u6709\"->\"\u5305\u62ec\",\n//     \"\u65b9\u6cd5\"->\"\u6cd5\" for the two trailing mentions).\n//  2. \"\u4f7f\u7528\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\u5bfb\u627e\u51cf\u5c11\u586b\u5165\u6392\u5217\" -> \"\u7528\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u6cd5\u8ba1\u7b97\u51cf\u5c11\u586b\u5165\u6392\u5217\"\n//  3. \"\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff0c\u9009\u62e9\u66f4\u4f18\u7684\u7ed3\u679c\u8fd4\u56de\" -> \"\u5d4c\u5957\u5206\u5272\u6cd5\uff0c\u9009\u62e9\u66f4\u4f18\u7684\u7ed3\u679c\u8fd4\u56de\"\n//  4. \"\u5f97\u5230\u7684\u5206\u6790\u7ed3\u679c\uff08\u6d88\u53bb\u6811\u53ca\u5176\u884d\u751f\u4fe1\u606f\uff09\u76f4\u63a5\u6307\u5bfc\" ->\n//     \"\u5f97\u5230\u7684\u5206\u6790\u7ed3\u679c\uff0c\u5373\u6d88\u53bb\u6811\u53ca\u5176\u884d\u751f\u4fe1\u606f\u5c06\u76f4\u63a5\u6307\u5bfc\"\n//  5. \"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\u7f57\u5217\u5982\" -> \"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\" + footnote + \"\u7f57\u5217\u5982\"\n\nconst body = context.document.body;\n\n// --- 1. CHOLMOD ordering-method sentence -------------------------------\nconst r1 = body.search(\n  \"\u4e3b\u8981\u6709\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\uff08CHOLMOD_AMD)\uff0c\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff08CHOLMOD_METIS\uff0c\" +\n    \"CHOLMOD_NESDIS\uff09\u7b49\u3002CHOLMOD\u7684\u9ed8\u8ba4\u5148\u5c1d\u8bd5\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\uff0c\u5982\u679c\u6548\u679c\u4e0d\u4f73\uff0c\" +\n    \"\u518d\u5c1d\u8bd5\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff0c\u8be6\u89c1\u7b97\u6cd5\",\n  { matchCase: true }\n);\nr1.load(\"text\");\nawait context.sync();\nif (r1.items.length !== 1) {\n  throw new Error(\"expected exactly one match for hunk 1, got \" + r1.items.length);\n}\nr1.items[0].insertText(\n  \"\u5305\u62ec\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\uff08Approximate minimum degree\uff09\uff0c\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff08\" +\n    \"Nested dissection\uff09\u7b49\u3002CHOLMOD\u7684\u9ed8\u8ba4\u5148\u5c1d\u8bd5\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u6cd5\uff0c\u5982\u679c\u6548\u679c\u4e0d\u4f73\" +\n    \"\u518d\u5c1d\u8bd5\u5d4c\u5957\u5206\u5272\u6cd5\uff0c\u8be6\u89c1\u7b97\u6cd5\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 2. \"\u4f7f\u7528\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\u5bfb\u627e\u51cf\u5c11\u586b\u5165\u6392\u5217\" ----------------------------\nconst r2 = body.search(\"\u4f7f\u7528\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\u5bfb\u627e\u51cf\u5c11\u586b\u5165\u6392\u5217\", { matchCase: true });\nr2.load(\"text\");\nawait context.sync();\nif (r2.items.length !== 1) {\n  throw new Error(\"expected exactly one match for hunk 2, got \" + r2.items.length);\n}\nr2.items[0].insertText(\"\u7528\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u6cd5\u8ba1\u7b97\u51cf\u5c11\u586b\u5165\u6392\u5217\", \"Replace\");\nawait context.sync();\n\n// --- 3. \"\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff0c\u9009\u62e9\u66f4\u4f18\u7684\u7ed3\u679c\u8fd4\u56de\" ---------------------------------\nconst r3 = body.search(\"\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff0c\u9009\u62e9\u66f4\u4f18\u7684\u7ed3\u679c\u8fd4\u56de\", { matchCase: true });\nr3.load(\"text\");\nawait context.sync();\nif (r3.items.length !== 1) {\n  throw new Error(\"expected exactly one match for hunk 3, got \" + r3.items.length);\n}\nr3.items[0].insertText(\"\u5d4c\u5957\u5206\u5272\u6cd5\uff0c\u9009\u62e9\u66f4\u4f18\u7684\u7ed3\u679c\u8fd4\u56de\", \"Replace\");\nawait context.sync();\n\n// --- 4. \"\u5f97\u5230\u7684\u5206\u6790\u7ed3\u679c\uff08\u6d88\u53bb\u6811\u53ca\u5176\u884d\u751f\u4fe1\u606f\uff09\u76f4\u63a5\u6307\u5bfc\" ----------------------\nconst r4 = body.search(\"\u5f97\u5230\u7684\u5206\u6790\u7ed3\u679c\uff08\u6d88\u53bb\u6811\u53ca\u5176\u884d\u751f\u4fe1\u606f\uff09\u76f4\u63a5\u6307\u5bfc\", { matchCase: true });\nr4.load(\"text\");\nawait context.sync();\nif (r4.items.length !== 1) {\n  throw new Error(\"expected exactly one match for hunk 4, got \" + r4.items.length);\n}\nr4.items[0].insertText(\"\u5f97\u5230\u7684\u5206\u6790\u7ed3\u679c\uff0c\u5373\u6d88\u53bb\u6811\u53ca\u5176\u884d\u751f\u4fe1\u606f\u5c06\u76f4\u63a5\u6307\u5bfc\", \"Replace\");\nawait context.sync();\n\n// --- 5. \"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\u7f57\u5217\u5982\" -> split + footnote ------------------------\n// The phrase is unique in the document; search for the shorter prefix\n// \"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\" scoped to the paragraph that contains the full phrase so\n// we get a collapsed insertion point right before \"\u7f57\u5217\u5982\".\nconst r5 = body.search(\"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\u7f57\u5217\u5982\", { matchCase: true });\nr5.load(\"text\");\nawait context.sync();\nif (r5.items.length !== 1) {\n  throw new Error(\"expected exactly one match for hunk 5, got \" + r5.items.length);\n}\nconst para = r5.items[0].paragraphs.getFirst();\nconst scoped = para.search(\"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\", { matchCase: true });\nscoped.load(\"text\");\nawait context.sync();\nif (scoped.items.length !== 1) {\n  throw new Error(\"expected exactly one scoped match for hunk 5, got \" + scoped.items.length);\n}\nconst dataRange = scoped.items[0];\nconst insertionPoint = dataRange.getRange(\"End\");\nconst footnote = insertionPoint.insertFootnote(\"\u6570\u636e\u6765\u6e90\u4e8e https://sparse.tamu.edu/\");\nawait context.sync();\n\n// Best-effort: align the new footnote's paragraph/character style with the\n// existing footnote style used elsewhere in the document (\"footnote text\" /\n// \"footnote reference\" map to this document's custom a8 / a9 style ids).\nfootnote.body.paragraphs.load(\"style\");\nawait context.sync();\nconst fnPara = footnote.body.paragraphs.getFirst();\nfnPara.style = \"footnote text\";\nfnPara.font.size = 8;\nawait context.sync();\n", "ps1": "# Apply the report.docx content edits described by the diff:\n#  1. Rewrite the CHOLMOD ordering-method sentence (expand CHOLMOD_AMD /\n#     CHOLMOD_METIS into spelled-out English names, \"\u4e3b\u8981\u6709\"->\"\u5305\u62ec\",\n#     \"\u65b9\u6cd5\"->\"\u6cd5\" for the two trailing mentions).\n#  2. \"\u4f7f\u7528\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\u5bfb\u627e\u51cf\u5c11\u586b\u5165\u6392\u5217\" -> \"\u7528\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u6cd5\u8ba1\u7b97\u51cf\u5c11\u586b\u5165\u6392\u5217\"\n#  3. \"\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff0c\u9009\u62e9\u66f4\u4f18\u7684\u7ed3\u679c\u8fd4\u56de\" -> \"\u5d4c\u5957\u5206\u5272\u6cd5\uff0c\u9009\u62e9\u66f4\u4f18\u7684\u7ed3\u679c\u8fd4\u56de\"\n#  4. \"\u5f97\u5230\u7684\u5206\u6790\u7ed3\u679c\uff08\u6d88\u53bb\u6811\u53ca\u5176\u884d\u751f\u4fe1\u606f\uff09\u76f4\u63a5\u6307\u5bfc\" ->\n#     \"\u5f97\u5230\u7684\u5206\u6790\u7ed3\u679c\uff0c\u5373\u6d88\u53bb\u6811\u53ca\u5176\u884d\u751f\u4fe1\u606f\u5c06\u76f4\u63a5\u6307\u5bfc\"\n#  5. \"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\u7f57\u5217\u5982\" -> \"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\" + footnote + \"\u7f57\u5217\u5982\"\n\n$d = $word.ActiveDocument\n\n# --- 1. CHOLMOD ordering-method sentence --------------------------------\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\n    \"\u4e3b\u8981\u6709\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\uff08CHOLMOD_AMD)\uff0c\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff08CHOLMOD_METIS\uff0c\" +\n    \"CHOLMOD_NESDIS\uff09\u7b49\u3002CHOLMOD\u7684\u9ed8\u8ba4\u5148\u5c1d\u8bd5\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\uff0c\u5982\u679c\u6548\u679c\u4e0d\u4f73\uff0c\" +\n    \"\u518d\u5c1d\u8bd5\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff0c\u8be6\u89c1\u7b97\u6cd5\")\nif (-not $found1) {\n    throw \"hunk 1 search text not found\"\n}\n$rng1.Text = (\n    \"\u5305\u62ec\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\uff08Approximate minimum degree\uff09\uff0c\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff08\" +\n    \"Nested dissection\uff09\u7b49\u3002CHOLMOD\u7684\u9ed8\u8ba4\u5148\u5c1d\u8bd5\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u6cd5\uff0c\u5982\u679c\u6548\u679c\u4e0d\u4f73\" +\n    \"\u518d\u5c1d\u8bd5\u5d4c\u5957\u5206\u5272\u6cd5\uff0c\u8be6\u89c1\u7b97\u6cd5\")\n\n# --- 2. \"\u4f7f\u7528\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\u5bfb\u627e\u51cf\u5c11\u586b\u5165\u6392\u5217\" -----------------------------\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"\u4f7f\u7528\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u65b9\u6cd5\u5bfb\u627e\u51cf\u5c11\u586b\u5165\u6392\u5217\")\nif (-not $found2) {\n    throw \"hunk 2 search text not found\"\n}\n$rng2.Text = \"\u7528\u8fd1\u4f3c\u6700\u5c0f\u5ea6\u6570\u6cd5\u8ba1\u7b97\u51cf\u5c11\u586b\u5165\u6392\u5217\"\n\n# --- 3. \"\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff0c\u9009\u62e9\u66f4\u4f18\u7684\u7ed3\u679c\u8fd4\u56de\" ----------------------------------\n$rng3 = $d.Content\n$found3 = $rng3.Find.Execute(\"\u5d4c\u5957\u5206\u5272\u65b9\u6cd5\uff0c\u9009\u62e9\u66f4\u4f18\u7684\u7ed3\u679c\u8fd4\u56de\")\nif (-not $found3) {\n    throw \"hunk 3 search text not found\"\n}\n$rng3.Text = \"\u5d4c\u5957\u5206\u5272\u6cd5\uff0c\u9009\u62e9\u66f4\u4f18\u7684\u7ed3\u679c\u8fd4\u56de\"\n\n# --- 4. \"\u5f97\u5230\u7684\u5206\u6790\u7ed3\u679c\uff08\u6d88\u53bb\u6811\u53ca\u5176\u884d\u751f\u4fe1\u606f\uff09\u76f4\u63a5\u6307\u5bfc\" -----------------------\n$rng4 = $d.Content\n$found4 = $rng4.Find.Execute(\"\u5f97\u5230\u7684\u5206\u6790\u7ed3\u679c\uff08\u6d88\u53bb\u6811\u53ca\u5176\u884d\u751f\u4fe1\u606f\uff09\u76f4\u63a5\u6307\u5bfc\")\nif (-not $found4) {\n    throw \"hunk 4 search text not found\"\n}\n$rng4.Text = \"\u5f97\u5230\u7684\u5206\u6790\u7ed3\u679c\uff0c\u5373\u6d88\u53bb\u6811\u53ca\u5176\u884d\u751f\u4fe1\u606f\u5c06\u76f4\u63a5\u6307\u5bfc\"\n\n# --- 5. \"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\u7f57\u5217\u5982\" -> split + footnote -------------------------\n$rng5 = $d.Content\n$found5 = $rng5.Find.Execute(\"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\")\nif (-not $found5) {\n    throw \"hunk 5 search text not found\"\n}\n$rng5.Collapse(0)   # wdCollapseEnd -> collapse to right after \"\u5b9e\u9a8c\u4e2d\u7528\u5230\u7684\u6570\u636e\"\n$fn = $d.Footnotes.Add($rng5, \"\", \"\u6570\u636e\u6765\u6e90\u4e8e https://sparse.tamu.edu/\")\n\n# Best-effort: align the new footnote's paragraph/character style with the\n# existing footnote style used elsewhere in the document (\"footnote text\" /\n# \"footnote reference\" map to this document's custom a8 / a9 style ids).\n$fn.Range.set_Style(\"footnote text\")\n$fn.Range.Font.Size = 8\n"}
